$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1115
$ws1.Range("F4").Value = 1791
$ws1.Range("F5").Value = 791
$ws1.Range("F6").Value = 324
$ws1.Range("F7").Value = 210

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1115
$ws4.Range("F4").Value = 1791
$ws4.Range("F6").Value = 791
$ws4.Range("F7").Value = 324
$ws4.Range("F8").Value = 210
